# Apply updates produced by refreshing CalcExcessMortYLL after changes to its
# CalcExcessMortality dependency. This touches the "raw" numeric columns
# (H, O, P, Q, R, W, X) as well as the derived "pretty" text columns
# (T = Excess Mortality ±, U = P_score ±, Y = Excess_mortality_per_10^5 ±)
# on the Female, Male and Total sheets for Belgium (row 3), Czechia (row 7)
# and Slovakia (row 28).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Female sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Female")

# Belgium (row 3)
$ws.Range("H3").Value = 55279
$ws.Range("O3").Value = 9992.200000000001
$ws.Range("T3").Value = "9992.2 (±615.7)"
$ws.Range("W3").Value = 171.1
$ws.Range("X3").Value = 10.5
$ws.Range("Y3").Value = "171.1(±10.5)"

# Czechia (row 7)
$ws.Range("H7").Value = 53452
$ws.Range("O7").Value = 8728.6
$ws.Range("T7").Value = "8728.6 (±750.6)"

# Slovakia (row 28)
$ws.Range("H28").Value = 24290
$ws.Range("O28").Value = 3012
$ws.Range("T28").Value = "3012.0 (±308.0)"
$ws.Range("Q28").Value = 14.2
$ws.Range("R28").Value = 1.7
$ws.Range("U28").Value = "14.2% (±1.7%)"
$ws.Range("W28").Value = 107.9
$ws.Range("Y28").Value = "107.9(±11.0)"

# ---------------------------------------------------------------------
# Male sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Male")

# Belgium (row 3)
$ws.Range("H3").Value = 53114
$ws.Range("O3").Value = 9325
$ws.Range("T3").Value = "9325.0 (±335.8)"

# Czechia (row 7)
$ws.Range("H7").Value = 57009
$ws.Range("O7").Value = 10616.2
$ws.Range("T7").Value = "10616.2 (±643.0)"
$ws.Range("X7").Value = 12.2
$ws.Range("Y7").Value = "201.4(±12.2)"

# Slovakia (row 28)
$ws.Range("H28").Value = 25867
$ws.Range("O28").Value = 3422
$ws.Range("T28").Value = "3422.0 (±257.9)"
$ws.Range("W28").Value = 128.4
$ws.Range("Y28").Value = "128.4(±9.7)"

# ---------------------------------------------------------------------
# Total sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Total")

# Belgium (row 3)
$ws.Range("H3").Value = 108393
$ws.Range("O3").Value = 19317.2
$ws.Range("T3").Value = "19317.2 (±924.0)"
$ws.Range("X3").Value = 8.1
$ws.Range("Y3").Value = "167.6(±8.1)"

# Slovakia (row 28)
$ws.Range("H28").Value = 50157
$ws.Range("O28").Value = 6434
$ws.Range("T28").Value = "6434.0 (±544.6)"
$ws.Range("W28").Value = 117.9
$ws.Range("X28").Value = 10
$ws.Range("Y28").Value = "117.9(±10.0)"
